$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.354.89"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "3.520.20"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.58%  "
$ws.Range("E9").Value = "  +6.75%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "4.127.25"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E14").Value = "  +4.10%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "67.247.13"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "3.495.65"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E38").Value = "  +7.61%  "
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("E43").Value = "  +5.78%  "
$ws.Range("D44").Value = "2.842.06"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +1.07%  "
$ws.Range("E51").Value = "  -0.70%  "

# Numeric-looking values in column D must stay as text (matching the
# original inlineStr cells), so force text storage then restore the
# default (unstyled) cell appearance.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.133"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.438"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("D15").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.892"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0756"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.53"
$ws.Range("D50").Style = "Normal"
